# Normalise the "Recorded By" (column G) attendance-log strings on the
# "Session Analysis Results" sheet: a handful of distinct author/system
# combinations need their comma-separated entries reordered.
#
# Mapping applied (old value -> new value):
#   "backup@backdoor.com, System, system" -> "system, System, backup@backdoor.com"
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System"         -> "System, backup@backdoor.com"
#
# Any other "Recorded By" value (e.g. "System", "System, admin@admin.com",
# "dnasr281@gmail.com", "dnasr281@gmail.com, admin@admin.com") is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$recordedByMap = @{
    "backup@backdoor.com, System, system" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = $cell.Text

    if ($recordedByMap.ContainsKey($current)) {
        $cell.Value = $recordedByMap[$current]
    }
}
